$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 236, pushing existing row 236 (and below) down to 237.
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new Haba price record.
$ws.Cells.Item(236, 1).Value = 3
$ws.Cells.Item(236, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 45215
$ws.Cells.Item(236, 5).Value = 5
$ws.Cells.Item(236, 6).Value = 100112026
$ws.Cells.Item(236, 7).Value = "Haba"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 40
$ws.Cells.Item(236, 11).Value = 12000
$ws.Cells.Item(236, 12).Value = 12000
$ws.Cells.Item(236, 13).Value = 12000
$ws.Cells.Item(236, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(236, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(236, 16).Value = 480
$ws.Cells.Item(236, 17).Value = 25
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(236, 4).NumberFormat = $ws.Cells.Item(237, 4).NumberFormat
